$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.992.35"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.778.80"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'585.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'161.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.33%  "
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.794.12"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "'0.399"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "3.268.34"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "'27.49"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "63.944.04"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("E17").Value = "  +6.70%  "
$ws.Range("D18").Value = "2.785.14"
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("D19").Value = "'12.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.44%  "
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("D21").Value = "'367.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'0.564"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.06%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'67.42"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'0.178"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.40%  "
$ws.Range("D27").Value = "'8.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "0.0₃0968"
$ws.Range("E28").Value = "  +14.74%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'2.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'7.36"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("D32").Value = "'1.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.79%  "
$ws.Range("D33").Value = "'172.52"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'5.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.58%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'20.87"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").Value = "'1.03"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'342.84"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").Value = "'6.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +11.48%  "
$ws.Range("D43").Value = "'39.95"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").Value = "'22.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.34%  "
$ws.Range("E45").Value = "  +5.70%  "
$ws.Range("E46").Value = "  +3.73%  "
$ws.Range("D47").Value = "'0.655"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").Value = "'138.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "2.176.20"
$ws.Range("E51").Value = "  +2.15%  "
